# "add the time management" - append new diary/time-log entries (rows 148-159)
# to Sheet1, mirroring the existing "date / time / place / note" layout used
# throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# xlPasteFormats constant used below to copy an existing cell's style
# (number format / font) onto a newly written cell without disturbing its
# value.
$xlPasteFormats = -4122

function Copy-Style($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy()
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats)
}

# ---- Row 148 : 十二月十七号 ----
$ws.Range("A148").Value = "十二月十七号"
$ws.Range("C148").Value = "在家编产品编号"

# ---- Row 149 : 十二月十八号 ----
$ws.Range("A149").Value = "十二月十八号"
Copy-Style "B147" "B149"
$ws.Range("B149").Value = 0.39583333333333331
$ws.Range("C149").Value = "见二师兄，一起学习，lentern,chrome,输入法，计算机概论，日语第八课"

# ---- Row 150 : 十二月十九号 ----
$ws.Range("A150").Value = "十二月十九号"
$ws.Range("C150").Value = "在家上传产品到国际站"

# ---- Row 151 : 十二月二十号 ----
$ws.Range("A151").Value = "十二月二十号"
Copy-Style "B147" "B151"
$ws.Range("B151").Value = 0.41666666666666669
$ws.Range("C151").Value = "北邮漫咖啡"
$ws.Range("D151").Value = "上传产品，定去葫芦岛行程"

# ---- Row 152 : 十二月二十一号 ----
$ws.Range("A152").Value = "十二月二十一号"
Copy-Style "B147" "B152"
$ws.Range("B152").Value = 0.41666666666666669
$ws.Range("C152").Value = "茶馆"
$ws.Range("D152").Value = "上传产品"

# ---- Row 153 : 十二月二十二号 ----
$ws.Range("A153").Value = "十二月二十二号"
Copy-Style "B147" "B153"
$ws.Range("B153").Value = 0.41666666666666669
$ws.Range("C153").Value = "北邮"
$ws.Range("D153").Value = "和哈立德一起出发"

# ---- Row 154 : (continuation of 十二月二十二号) ----
Copy-Style "B147" "B154"
$ws.Range("B154").Value = 0.80555555555555547
$ws.Range("C154").Value = "辽工大"
$ws.Range("D154").Value = "见院长浅谈只会矿山项目"

# ---- Row 155 : 十二月二十三号 ----
$ws.Range("A155").Value = "十二月二十三号"
Copy-Style "B147" "B155"
$ws.Range("B155").Value = 0.52083333333333337
$ws.Range("C155").Value = "古城"

# ---- Row 156 : —二十五号 ----
Copy-Style "C129" "A156"
$ws.Range("A156").Value = "—二十五号"
Copy-Style "B147" "B156"
$ws.Range("B156").Value = 0.58333333333333337
$ws.Range("C156").Value = "辽工大"
$ws.Range("D156").Value = "等院长，见了书记，老师，同学"

# ---- Row 157 : 个人工作时间记录 ----
Copy-Style "C129" "A157"
$ws.Range("A157").Value = "个人工作时间记录"

# ---- Row 158 : 十二月二十六号 ----
Copy-Style "C129" "A158"
$ws.Range("A158").Value = "十二月二十六号"
Copy-Style "B147" "B158"
$ws.Range("B158").Value = 0.4375
Copy-Style "C129" "C158"
$ws.Range("C158").Value = "茶馆"
Copy-Style "C129" "D158"
$ws.Range("D158").Value = "更改产品信息（40个产品）"

# ---- Row 159 : (continuation) ----
Copy-Style "B147" "B159"
$ws.Range("B159").Value = 0.5
Copy-Style "C129" "D159"
$ws.Range("D159").Value = "更改完"

# Column A needed to grow a bit wider to fit the longer date labels such as
# "十二月二十六号".
$ws.Range("A1").ColumnWidth = 16.75

# Page setup touched (paper size / orientation) as part of the edit.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection on the newly added summary row, like the recorded
# session did.
$ws.Range("A157").Select()
